# Add cascade FE data function
# - Clear placeholder "FE function" flag/weight values that were temporarily
#   populated on the Compartments and Characteristics sheets, restoring the
#   blank cells (H2:H5 on Compartments, E2:E9 on Characteristics).
# - Update sheet selections / active tab to reflect the new working sheet
#   (Compartments) instead of Cascades.

$wb = $excel.ActiveWorkbook

$wsCompartments = $wb.Worksheets.Item("Compartments")
$wsCharacteristics = $wb.Worksheets.Item("Characteristics")
$wsCascades = $wb.Worksheets.Item("Cascades")

# Clear the H2:H5 values on Compartments (keep formatting/style)
$wsCompartments.Range("H2:H5").ClearContents()

# Clear the E2:E9 values on Characteristics (keep formatting/style)
$wsCharacteristics.Range("E2:E9").ClearContents()

# Update selections
$wsCompartments.Range("L5").Select()
$wsCharacteristics.Range("E2:E9").Select()
$wsCascades.Range("D9").Select()

# Make Compartments the active/selected sheet (tabSelected) and
# deactivate Cascades' tab selection.
$wsCompartments.Activate()
$wsCompartments.Select()
